# June community call presentation.
#
# 1) Slide 13 ("Office 365 transition plan status") title: the title had
#    been typed as two runs ("Office 365 " + "transition plan status")
#    that collapse into a single run once the slide is re-saved.
# 2) Slide 13 body text box (inside "Group 32"): "the " and "Office 365 "
#    were two separate runs that collapse into a single run "the Office 365 ".
# 3) Slide 45 footer/credits text box: "2014" -> "2015" in the copyright
#    line, typed as "2015 " inserted right after "(c) " and "2014 " deleted
#    from the following run, leaving "Microsoft " in place.

$p = $ppt.ActivePresentation

# --- Slide 13: title run merge -------------------------------------------
$s13 = $p.Slides.Item(13)

$title = $s13.Shapes.Item(1)
$titleRange = $title.TextFrame.TextRange
$titleWhole = $titleRange.Characters(1, $titleRange.Length)
$titleWhole.Text = "Office 365 transition plan status"

# --- Slide 13: "the " + "Office 365 " run merge inside the grouped text box
$grp = $s13.Shapes.Item(2)
$capTextBox = $grp.GroupItems.Item(1)
$capRange = $capTextBox.TextFrame.TextRange
$capText = $capRange.Text
$mergeStart = $capText.IndexOf("the Office 365 ") + 1
$capSel = $capRange.Characters($mergeStart, 15)
$capSel.Text = "the Office 365 "

# --- Slide 45: copyright year 2014 -> 2015 --------------------------------
$s45 = $p.Slides.Item(45)
$credits = $s45.Shapes.Item(1)
$credRange = $credits.TextFrame.TextRange
$credText = $credRange.Text
$yearStart = $credText.IndexOf("2014 ") + 1
$yearSel = $credRange.Characters($yearStart, 5)
$yearSel.Text = "2015 "
